$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): Q_Responsibility (X) and Q_Responsibilitycode (Y) ---
$ws.Range("X1").Value = "Q_Responsibility"
$ws.Range("Y1").Value = "Q_Responsibilitycode"

# --- Ensure the Q_Responsibilitycode column (Y) stores its numeric-looking
#     codes as text, consistent with the other "...code" columns (e.g. G, W) ---
$ws.Range("Y2:Y79").NumberFormat = "@"

# --- Data values for rows 2-79 ---
$responsibility = @{
    2 = "Public authorities are completely responsible for flood protection"
    3 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    4 = "Public authorities are completely responsible for flood protection"
    5 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    6 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    7 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    8 = "Public authorities are completely responsible for flood protection"
    9 = "Public authorities are completely responsible for flood protection"
    10 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    11 = "Public authorities are completely responsible for flood protection"
    12 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    13 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    14 = "Public authorities are completely responsible for flood protection"
    15 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    16 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    17 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    18 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    19 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    20 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    21 = "Public authorities are completely responsible for flood protection"
    22 = "Public authorities are completely responsible for flood protection"
    23 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    24 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    25 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    26 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    27 = "Public authorities are completely responsible for flood protection"
    28 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    29 = "Public authorities are completely responsible for flood protection"
    30 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    31 = "Public authorities are completely responsible for flood protection"
    32 = "Public authorities are completely responsible for flood protection"
    33 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    34 = "Public authorities are completely responsible for flood protection"
    35 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    36 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    37 = "Public authorities are completely responsible for flood protection"
    38 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    39 = "Public authorities are completely responsible for flood protection"
    40 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    41 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    42 = "Public authorities are completely responsible for flood protection"
    43 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    44 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    45 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    46 = "Public authorities are completely responsible for flood protection"
    47 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    48 = "Public authorities are completely responsible for flood protection"
    49 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    50 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    51 = "Public authorities are completely responsible for flood protection"
    52 = "Public authorities and citizens are equally responsible for flood protection"
    53 = "Public authorities and citizens are equally responsible for flood protection"
    54 = "Public authorities are completely responsible for flood protection"
    55 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    56 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    57 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    58 = "Public authorities are completely responsible for flood protection"
    59 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    60 = "Public authorities are completely responsible for flood protection"
    61 = "Public authorities are completely responsible for flood protection"
    62 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    63 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    64 = "Public authorities and citizens are equally responsible for flood protection"
    65 = "Public authorities are completely responsible for flood protection"
    66 = "Public authorities are completely responsible for flood protection"
    67 = "Public authorities are completely responsible for flood protection"
    68 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    69 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    70 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    71 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    72 = "Public authorities are completely responsible for flood protection"
    73 = "Public authorities and citizens are equally responsible for flood protection"
    74 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    75 = "Public authorities are completely responsible for flood protection"
    76 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    77 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    78 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
    79 = "Public authorities are responsible and citizens somewhat responsible for flood protection"
}

$responsibilitycode = @{
    2 = "1"
    3 = "2"
    4 = "1"
    5 = "2"
    6 = "2"
    7 = "2"
    8 = "1"
    9 = "1"
    10 = "2"
    11 = "1"
    12 = "2"
    13 = "2"
    14 = "1"
    15 = "2"
    16 = "2"
    17 = "2"
    18 = "2"
    19 = "2"
    20 = "2"
    21 = "1"
    22 = "1"
    23 = "2"
    24 = "2"
    25 = "2"
    26 = "2"
    27 = "1"
    28 = "2"
    29 = "1"
    30 = "2"
    31 = "1"
    32 = "1"
    33 = "2"
    34 = "1"
    35 = "2"
    36 = "2"
    37 = "1"
    38 = "2"
    39 = "1"
    40 = "2"
    41 = "2"
    42 = "1"
    43 = "2"
    44 = "2"
    45 = "2"
    46 = "1"
    47 = "2"
    48 = "1"
    49 = "2"
    50 = "2"
    51 = "1"
    52 = "3"
    53 = "3"
    54 = "1"
    55 = "2"
    56 = "2"
    57 = "2"
    58 = "1"
    59 = "2"
    60 = "1"
    61 = "1"
    62 = "2"
    63 = "2"
    64 = "3"
    65 = "1"
    66 = "1"
    67 = "1"
    68 = "2"
    69 = "2"
    70 = "2"
    71 = "2"
    72 = "1"
    73 = "3"
    74 = "2"
    75 = "1"
    76 = "2"
    77 = "2"
    78 = "2"
    79 = "2"
}

for ($r = 2; $r -le 79; $r++) {
    $ws.Cells.Item($r, 24).Value = $responsibility[$r]
    $ws.Cells.Item($r, 25).Value = $responsibilitycode[$r]
}

Write-Host ("Updated X1:Y79 on sheet " + $ws.Name)